$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.875.74'
$ws.Range("E2").Value = '  +0.75%  '

$ws.Range("D3").Value = '1.642.54'
$ws.Range("E3").Value = '  +0.71%  '

$ws.Range("E4").Value = '  -0.67%  '

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '216.75'
$cell.Style = $origStyle
$ws.Range("E5").Value = '  -0.44%  '

$ws.Range("E6").Value = '  +2.12%  '

$ws.Range("E7").Value = '  -0.63%  '

$ws.Range("E8").Value = '  +2.09%  '

$ws.Range("E9").Value = '  +0.45%  '

$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '19.84'
$cell.Style = $origStyle
$ws.Range("E10").Value = '  +4.73%  '

$ws.Range("E11").Value = '  +0.36%  '

$ws.Range("D12").Value = '1.872.00'
$ws.Range("E12").Value = '  +0.61%  '

$ws.Range("D13").Value = '1.637.02'
$ws.Range("E13").Value = '  +0.16%  '

$ws.Range("E14").Value = '  +0.78%  '

$ws.Range("E15").Value = '  +1.68%  '

$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '66.45'
$cell.Style = $origStyle
$ws.Range("E16").Value = '  +3.99%  '

$ws.Range("D17").Value = '26.883.87'
$ws.Range("E17").Value = '  +0.78%  '

$ws.Range("D18").Value = '0.0₃0729'
$ws.Range("E18").Value = '  +1.16%  '

$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '219.48'
$cell.Style = $origStyle
$ws.Range("E19").Value = '  +3.92%  '

$ws.Range("E20").Value = '  -0.57%  '

$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.65'
$cell.Style = $origStyle
$ws.Range("E21").Value = '  +7.91%  '

$ws.Range("E22").Value = '  +2.14%  '

$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.43'
$cell.Style = $origStyle
$ws.Range("E23").Value = '  +3.86%  '

$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '9.19'
$cell.Style = $origStyle
$ws.Range("E24").Value = '  +0.59%  '

$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '145.88'
$cell.Style = $origStyle
$ws.Range("E25").Value = '  -0.40%  '

$ws.Range("E26").Value = '  -0.70%  '

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.40'
$cell.Style = $origStyle
$ws.Range("E27").Value = '  +5.80%  '

$ws.Range("E28").Value = '  +1.67%  '

$ws.Range("E29").Value = '  +2.19%  '

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0509'
$cell.Style = $origStyle
$ws.Range("E30").Value = '  +1.43%  '

$ws.Range("E31").Value = '  -0.35%  '

$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.35'
$cell.Style = $origStyle
$ws.Range("E32").Value = '  -0.28%  '

$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.00'
$cell.Style = $origStyle
$ws.Range("E33").Value = '  +2.29%  '

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.56'
$cell.Style = $origStyle
$ws.Range("E34").Value = '  +3.16%  '

$ws.Range("E35").Value = '  +0.12%  '

$ws.Range("D36").Value = '1.245.72'
$ws.Range("E36").Value = '  -0.96%  '

$ws.Range("E37").Value = '  +1.08%  '

$ws.Range("E38").Value = '  +3.76%  '

$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.833'
$cell.Style = $origStyle
$ws.Range("E39").Value = '  +4.10%  '

$ws.Range("E40").Value = '  -0.54%  '

$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.808'
$cell.Style = $origStyle
$ws.Range("E41").Value = '  +1.62%  '

$ws.Range("E42").Value = '  +2.42%  '

$ws.Range("D43").Value = '1.783.44'
$ws.Range("E43").Value = '  +0.65%  '

$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.08'
$cell.Style = $origStyle
$ws.Range("E44").Value = '  -3.41%  '

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '60.81'
$cell.Style = $origStyle
$ws.Range("E45").Value = '  +1.82%  '

$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '91.54'
$cell.Style = $origStyle
$ws.Range("E46").Value = '  +0.76%  '

$ws.Range("E47").Value = '  +0.94%  '

$ws.Range("D48").Value = '0.0₆0106'
$ws.Range("E48").Value = '  +18.48%  '

$ws.Range("E49").Value = '  -0.22%  '

$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0976'
$cell.Style = $origStyle
$ws.Range("E50").Value = '  +2.36%  '

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.59'
$cell.Style = $origStyle
$ws.Range("E51").Value = '  +2.04%  '
